$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date-like ("2024-01-08") and week ("01") strings auto-convert to numbers
# when assigned via .Value, so force them through a Text format then strip
# the resulting style so the cell matches the rest of the sheet (no s=).
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = "2024-01-08"
$ws.Range("A36").ClearFormats()

$ws.Range("B36").Value = "22:15:08"
$ws.Range("C36").Value = "Monday"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "01"
$ws.Range("D36").ClearFormats()

$ws.Range("E36").Value = 139542
$ws.Range("F36").Value = 142941
$ws.Range("G36").Value = 172486
$ws.Range("H36").Value = 147207
$ws.Range("I36").Value = -1
$ws.Range("J36").Value = 118259
$ws.Range("K36").Value = 224721
$ws.Range("L36").Value = 249917
$ws.Range("M36").Value = 185080
$ws.Range("N36").Value = 110390
$ws.Range("O36").Value = 40667
$ws.Range("P36").Value = 30821
$ws.Range("Q36").Value = 72430
$ws.Range("R36").Value = -1
$ws.Range("S36").Value = 42188
$ws.Range("T36").Value = -1
